# Apply manual quantity updates (MAJ manuelle) to FOURNISSEUR_A sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - P001: 50 -> 8
$ws.Range("B2").Value = 8
$ws.Range("C2").Value = 8
$ws.Range("D2").Value = 8

# Row 3 - P002: 0 -> 4
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 4

# Row 7 - P006: 0 -> 10
$ws.Range("B7").Value = 10
$ws.Range("C7").Value = 10
$ws.Range("D7").Value = 10
